$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same rows of data (F2, F3, F6)
# that need their "想去人数" (want-to-go count) incremented by 1.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2956
    $ws.Range("F3").Value = 733
    $ws.Range("F6").Value = 1694
}
